$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings in column D stay as text (matches source formatting)
$textCells = @("D4","D5","D6","D7","D8","D10","D12","D16","D19","D20","D22","D23","D24","D25","D26","D27","D28","D30","D32","D35","D36","D40","D41","D42","D43","D44","D45","D46","D47","D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '57.919.41'
$ws.Range("E2").Value = '  -1.96%  '
$ws.Range("D3").Value = '2.456.07'
$ws.Range("E3").Value = '  -2.02%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '516.73'
$ws.Range("E5").Value = '  -3.70%  '
$ws.Range("D6").Value = '132.29'
$ws.Range("E6").Value = '  -2.73%  '
$ws.Range("D7").Value = '0.998'
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = '0.556'
$ws.Range("E8").Value = '  -1.79%  '
$ws.Range("D9").Value = '2.464.48'
$ws.Range("E9").Value = '  -2.12%  '
$ws.Range("D10").Value = '0.0982'
$ws.Range("E10").Value = '  -2.85%  '
$ws.Range("E11").Value = '  -0.08%  '
$ws.Range("D12").Value = '5.26'
$ws.Range("E12").Value = '  -1.12%  '
$ws.Range("E13").Value = '  -2.30%  '
$ws.Range("D14").Value = '2.892.29'
$ws.Range("E14").Value = '  -2.00%  '
$ws.Range("D15").Value = '57.849.68'
$ws.Range("E15").Value = '  -1.75%  '
$ws.Range("D16").Value = '22.22'
$ws.Range("E16").Value = '  -3.28%  '
$ws.Range("D18").Value = '2.451.61'
$ws.Range("E18").Value = '  -2.46%  '
$ws.Range("D19").Value = '10.65'
$ws.Range("E19").Value = '  -3.61%  '
$ws.Range("D20").Value = '320.07'
$ws.Range("E20").Value = '  -0.90%  '
$ws.Range("E21").Value = '  -2.38%  '
$ws.Range("D22").Value = '0.999'
$ws.Range("E22").Value = '  +0.02%  '
$ws.Range("D23").Value = '5.73'
$ws.Range("E23").Value = '  -3.76%  '
$ws.Range("D24").Value = '64.33'
$ws.Range("E24").Value = '  -1.08%  '
$ws.Range("D25").Value = '0.408'
$ws.Range("E25").Value = '  -2.64%  '
$ws.Range("D26").Value = '0.997'
$ws.Range("E26").Value = '  +0.10%  '
$ws.Range("D27").Value = '0.160'
$ws.Range("E27").Value = '  -2.51%  '
$ws.Range("D28").Value = '7.32'
$ws.Range("E28").Value = '  -2.48%  '
$ws.Range("D29").Value = '0.0₃0738'
$ws.Range("E29").Value = '  -3.71%  '
$ws.Range("D30").Value = '166.50'
$ws.Range("E30").Value = '  -2.31%  '
$ws.Range("E31").Value = '  -4.24%  '
$ws.Range("D32").Value = '6.21'
$ws.Range("E32").Value = '  -6.34%  '
$ws.Range("E33").Value = '  -1.00%  '
$ws.Range("E34").Value = '  +0.08%  '
$ws.Range("D35").Value = '0.999'
$ws.Range("D36").Value = '18.05'
$ws.Range("E36").Value = '  -1.50%  '
$ws.Range("E37").Value = '  -6.21%  '
$ws.Range("E38").Value = '  -2.68%  '
$ws.Range("E39").Value = '  -4.23%  '
$ws.Range("D40").Value = '36.20'
$ws.Range("E40").Value = '  -1.86%  '
$ws.Range("D41").Value = '0.788'
$ws.Range("E41").Value = '  -3.00%  '
$ws.Range("D42").Value = '3.43'
$ws.Range("E42").Value = '  -4.27%  '
$ws.Range("D43").Value = '270.72'
$ws.Range("E43").Value = '  -4.76%  '
$ws.Range("D44").Value = '4.93'
$ws.Range("E44").Value = '  -4.41%  '
$ws.Range("D45").Value = '0.590'
$ws.Range("E45").Value = '  -2.56%  '
$ws.Range("D46").Value = '124.88'
$ws.Range("E46").Value = '  -3.77%  '
$ws.Range("D47").Value = '0.0906'
$ws.Range("E47").Value = '  -1.80%  '
$ws.Range("E48").Value = '  -3.55%  '
$ws.Range("E49").Value = '  -4.22%  '
$ws.Range("D50").Value = '16.68'
$ws.Range("E50").Value = '  -3.96%  '
$ws.Range("D51").Value = '1.723.13'
$ws.Range("E51").Value = '  -1.78%  '
